# Latest updated input values
# Update the GridCapacityInvestmentCost value (column L, row 2) for Benin specs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Update the input value: GridCapacityInvestmentCost 0.1065 -> 0.19
$ws.Range("L2").Value = 0.19

# Reflect the author's final cursor/view position on the sheet
$ws.Range("L2").Select()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
